# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across specific resume
# bullet points, matching the target OOXML diff.

$d = $word.ActiveDocument

# Highlight color used throughout the diff: w:color w:val="2C3E50"
# Word's Font.Color takes a BGR-packed integer (standard VBA RGB() order),
# which serializes back out as the RRGGBB hex string "2C3E50".
$highlightColor = 5258796
$plusMinus = [char]0x00B1

function Find-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

function Highlight-Metrics($paragraph, [string[]]$substrings) {
    if ($null -eq $paragraph) {
        Write-Output "WARNING: paragraph not found; skipping highlight batch"
        return
    }
    $paraEnd = $paragraph.Range.End
    $cursor = $paragraph.Range.Start
    foreach ($sub in $substrings) {
        $rng = $d.Range($cursor, $paraEnd)
        $found = $rng.Find.Execute($sub, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = 1
            $rng.Font.Color = $highlightColor
            $cursor = $rng.End
        }
    }
}

# 1. "• Discovered systematic race coding errors ... from 23% to 64%"
$p1 = Find-ParagraphByText "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
Highlight-Metrics $p1 @("23%", "64%")

# 2. "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
$p2 = Find-ParagraphByText ("Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $plusMinus + "4.2% to " + $plusMinus + "2.1%")
Highlight-Metrics $p2 @("87%", "71%", ($plusMinus + "4.2%"), ($plusMinus + "2.1%"))

# 3. "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
$p3 = Find-ParagraphByText "Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Highlight-Metrics $p3 @("1,200")

# 4. "• Created comprehensive meta-analysis framework ... $400M Polling Consortium Database ... now valued at $1B+"
$p4 = Find-ParagraphByText "Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+"
Highlight-Metrics $p4 @("$400M", "$1B")

# 5. "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$p5 = Find-ParagraphByText "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Highlight-Metrics $p5 @("73.5%", "$4.7M")

# 6. "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (shorter variant,
#    distinguished from paragraph 2 above by the absence of the "reducing polling error..." tail)
$p6 = $null
foreach ($p in $d.Paragraphs) {
    if (($p.Range.Text -like "*Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%*") -and ($p.Range.Text -notlike "*reducing*")) {
        $p6 = $p
    }
}
Highlight-Metrics $p6 @("87%", "71%")

Write-Output "highlighting applied"
